# Natmi following Dr Hou advice
# Update Mif-Cxcr4 LR-pair stats: expressing-cell counts 1 -> 3 plus
# recomputed downstream totals/specificities for rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; F=1; G=28.401376;          H=85.204128;          I=0.2813463917610605;  J=0.2813463917610605;  K=3; L=1; M=224.2321046666667;  N=672.696314;   O=0.9009864013525987;  P=0.9009864013525988;  Q=6368.500315909355;  R=57316.50284318419;  S=0.2534892730463363;  T=0.2534892730463363 }
    3  = @{ E=3; F=1; G=28.401376;          H=85.204128;          I=0.2813463917610605;  J=0.2813463917610605;  K=3; L=1; M=1.023704333333333;   N=3.071113;     O=0.004113343558497904; P=0.004113343558497904; Q=29.07461168382933;  R=261.671505154464;   S=0.001157274368256986; T=0.001157274368256986 }
    4  = @{ E=3; F=1; G=28.401376;          H=85.204128;          I=0.2813463917610605;  J=0.2813463917610605;  K=3; L=1; M=23.61820766666667;   N=70.854623;    O=0.0949002550889034;  P=0.09490025508890343; Q=670.7895963870826;  R=6037.106367483744;  S=0.02669984434646719; T=0.0266998443464672 }
    5  = @{ E=3; F=1; G=14.83037466666667;  H=44.491124;          I=0.1469109243485705;  J=0.1469109243485705;  K=3; L=1; M=224.2321046666667;  N=672.696314;   O=0.9009864013525987;  P=0.9009864013525988;  Q=3325.446124501882;  R=29929.01512051693;  S=0.1323647450482024;  T=0.1323647450482024 }
    6  = @{ E=3; F=1; G=14.83037466666667;  H=44.491124;          I=0.1469109243485705;  J=0.1469109243485705;  K=3; L=1; M=1.023704333333333;   N=3.071113;     O=0.004113343558497904; P=0.004113343558497904; Q=15.18191881122356;  R=136.637269301012;   S=0.0006042951043421654; T=0.0006042951043421656 }
    7  = @{ E=3; F=1; G=14.83037466666667;  H=44.491124;          I=0.1469109243485705;  J=0.1469109243485705;  K=3; L=1; M=23.61820766666667;   N=70.854623;    O=0.0949002550889034;  P=0.09490025508890343; Q=350.2668686518057;  R=3152.401817866252;  S=0.01394188419602593; T=0.01394188419602594 }
    8  = @{ E=3; F=1; G=57.71632199999999;  H=173.148966;         I=0.571742683890369;   J=0.571742683890369;   K=3; L=1; M=224.2321046666667;  N=672.696314;   O=0.9009864013525987;  P=0.9009864013525988;  Q=12941.85235567904;  R=116476.6712011113;  S=0.5151323832580599;  T=0.51513238325806 }
    9  = @{ E=3; F=1; G=57.71632199999999;  H=173.148966;         I=0.571742683890369;   J=0.571742683890369;   K=3; L=1; M=1.023704333333333;   N=3.071113;     O=0.004113343558497904; P=0.004113343558497904; Q=59.08444893546199;  R=531.7600404191579;  S=0.002351774085898752; T=0.002351774085898753 }
    10 = @{ E=3; F=1; G=57.71632199999999;  H=173.148966;         I=0.571742683890369;   J=0.571742683890369;   K=3; L=1; M=23.61820766666667;   N=70.854623;    O=0.0949002550889034;  P=0.09490025508890343; Q=1363.156078752202;  R=12268.40470876982;  S=0.05425852654641028; T=0.05425852654641029 }
}

$columns = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $data.Keys) {
    $rowValues = $data[$row]
    foreach ($col in $columns) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}
